# Adds a new "Count" (int) field/column to the Shop table.
# The new field is inserted as column I, mirroring the shape of the
# existing meta-rows (Type/Public/Private/Save/Cache/Ref/Upload/Desc)
# and carrying a value of 1 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header
$ws.Range("I1").Value = "Count"

# Row 2: declared type
$ws.Range("I2").Value = "int"

# Rows 3-8: boolean meta flags (Public/Private/Save/Cache/Ref/Upload)
# Column I mirrors column H except for "Save" (row 5), which is TRUE.
$ws.Range("I3").Value = $false
$ws.Range("I4").Value = $false
$ws.Range("I5").Value = $true
$ws.Range("I6").Value = $false
$ws.Range("I7").Value = $false
$ws.Range("I8").Value = $false

# Row 9 (Desc) intentionally left blank - no description was authored.

# Copy column H's formatting for rows 1-8 onto column I so the new
# column visually matches the rest of the header/meta block.
$ws.Range("H1:H8").Copy()
$ws.Range("I1:I8").PasteSpecial(-4122) # xlPasteFormats

# Data rows 10-71: every shop entry gets a Count of 1.
for ($r = 10; $r -le 71; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
}

# Extend the boolean-list data validation (TRUE/FALSE) already applied
# to B7:H8 so it also covers the new column I for the same rows.
$ws.Range("I7:I8").Validation.Add(3, 1, 3, "TRUE,FALSE")
